$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B:H cells carry NumberFormat "0" (the integer format behind style
# index 1 in the original file). Apply the same format to the new column I
# (rows 1-24) so it resolves to the identical style.
$ws.Range("I1:I24").NumberFormat = "0"

# Add new row 25 with label "pasture_nr" and values for columns E:H
$ws.Range("A25").Value = "pasture_nr"

$ws.Range("B25:D25").NumberFormat = "0"

$ws.Range("E25").Value = 13.820106506347656
$ws.Range("F25").Value = 14.501720428466797
$ws.Range("G25").Value = 17.01146125793457
$ws.Range("H25").Value = 17.766544342041016
$ws.Range("E25:H25").NumberFormat = "0"

$ws.Range("I25").NumberFormat = "0"
